$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.306.20"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.928.81"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").Value = "4.550.99"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "3.923.04"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "68.371.99"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +17.03%  "
$ws.Range("E26").Value = "  +12.08%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "718.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +19.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "0.0₃0891"
$ws.Range("E36").Value = "  +13.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.418"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +23.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.149"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.55%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  +6.87%  "
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.95%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0355"
$ws.Range("E49").Value = "  +34.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("E51").Value = "  -0.60%  "
